# Update Pat Cummins match-by-match stats (runs, balls, fours, sixes)
# for rows 2-11 on the active sheet. Values in this sheet are stored as
# text (numbers-as-text), so we keep the cells formatted as Text before
# writing the new values to avoid Excel auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the "numbers stored as text" representation used by the rest
# of the sheet for this data block.
$ws.Range("C2:F12").NumberFormat = "@"

# row, runs, balls, fours, sixes
$updates = @(
    , @(2,  "5",  "4", "1", "0")
    , @(3,  "1",  "8", "0", "0")
    , @(4,  "5",  "4", "0", "0")
    , @(5,  "4", "17", "0", "0")
    , @(6, "17",  "9", "1", "1")
    , @(7, "33", "12", "1", "4")
    , @(8, "53", "36", "5", "2")
    , @(9,  "1",  "3", "0", "0")
    , @(10, "15", "11", "0", "1")
    , @(11, "0",  "0", "0", "0")
)

foreach ($entry in $updates) {
    $row = $entry[0]
    $ws.Range("C$row").Value = $entry[1]
    $ws.Range("D$row").Value = $entry[2]
    $ws.Range("E$row").Value = $entry[3]
    $ws.Range("F$row").Value = $entry[4]
}
